# THIS IS FOR PART 2 GITHUB DELIVERABLE.docx
#
# The chat exchange gets a reply: two blank lines are added after the
# existing sentence, followed by a new line "Ok I go now." The _GoBack
# bookmark (which Word uses to remember the last edited spot) is moved so
# that it spans from the very beginning of the text to the very end of the
# newly typed sentence, reflecting the fact that the whole passage - from
# the first keystroke to the last - was just edited.

$d = $word.ActiveDocument

# Find the end of the existing sentence ("...right now") and collapse the
# range to a single point right after it, so new text/paragraphs land there.
$findRange = $d.Content
$findRange.Find.Execute("right now") | Out-Null
$insertionPoint = $findRange.Duplicate
$insertionPoint.Collapse(0)   # wdCollapseEnd

# Press Enter twice (two empty paragraphs) then type the new sentence.
$insertionPoint.InsertAfter("`r`r`rOk I go now.")

# Move the _GoBack bookmark: delete the old (collapsed) one and re-create it
# so it starts at the very top of the document and ends at the very end,
# i.e. it now wraps everything that was touched during this editing pass.
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()
$fullRange = $d.Range(0, $d.Content.End)
$d.Bookmarks.Add("_GoBack", $fullRange)
